$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the shared string for row 22 first, so it becomes the earlier
# entry in the shared strings table (matches original authoring order).
$ws.Cells.Item(22, 4).Value = "Implementation of Round-Robin. GNU FDL added to manual"
$ws.Cells.Item(21, 4).Value = "Manual set up from LaTeX template"

# Row 21: 2012-10-08 (serial 41190), 1h effort, note about manual setup from LaTeX template
$ws.Cells.Item(21, 1).Value = 41190
$ws.Cells.Item(21, 2).Value = 1

# Row 22: 2012-10-10 (serial 41192), 2h effort, note about round-robin implementation
$ws.Cells.Item(22, 1).Value = 41192
$ws.Cells.Item(22, 2).Value = 2

# Reuse the same date style (column A style "1") already used by rows 2-20,
# rather than letting NumberFormat create a duplicate numFmt entry.
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A21:A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update selection to reflect the new active cell after data entry
$ws.Range("A23").Select()
